$wb = $excel.ActiveWorkbook

# --- Update existing sheet "Redis_Single_Filter": add row 4 ---
$ws1 = $wb.Worksheets.Item("Redis_Single_Filter")

$ws1.Cells.Item(4, 1).Value = "InsultFilter"
$ws1.Cells.Item(4, 2).Value = "Redis"
$ws1.Cells.Item(4, 3).Value = "Single-node"
$ws1.Cells.Item(4, 4).Value = 1
$ws1.Cells.Item(4, 5).Value = 1000
$ws1.Cells.Item(4, 6).Value = 3.23
$ws1.Cells.Item(4, 7).Value = 310.06

# --- Add new sheet "Redis_Single_Service" (positioned after the first sheet) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Redis_Single_Service"

$ws2.Cells.Item(1, 1).Value = "Test"
$ws2.Cells.Item(1, 2).Value = "Middleware"
$ws2.Cells.Item(1, 3).Value = "Mode"
$ws2.Cells.Item(1, 4).Value = "Clients"
$ws2.Cells.Item(1, 5).Value = "Num Tasks"
$ws2.Cells.Item(1, 6).Value = "Temps Total (s)"
$ws2.Cells.Item(1, 7).Value = "RPS"

# Reuse the same header formatting (bold, centered, bordered) as sheet1
$ws1.Range("A1:G1").Copy()
$ws2.Range("A1:G1").PasteSpecial(-4122)

$ws2.Cells.Item(2, 1).Value = "InsultService"
$ws2.Cells.Item(2, 2).Value = "Redis"
$ws2.Cells.Item(2, 3).Value = "Single-node"
$ws2.Cells.Item(2, 4).Value = 1
$ws2.Cells.Item(2, 5).Value = 1000
$ws2.Cells.Item(2, 6).Value = 1.05
$ws2.Cells.Item(2, 7).Value = 950.27
